$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the 'Pending' column header to 'Ready'
$ws.Range("M1").Value = "Ready"

# Update the data in the last row to reflect the rename
$ws.Range("D16").Value = "ready from db"
$ws.Range("M16").Value = "n"

# Move the active selection near the edited row (cosmetic: matches the edit)
$ws.Range("D17").Select()
